# Update NATMI LR-pair data (Efna5-Ephb2) with new TPM-based values.
# Adds a new "ECs" sending-cluster block (rows 2-4) and refreshes the
# previously existing FAPs/MuSCs sending-cluster rows (now rows 5-10)
# with recomputed specificity/weight figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs | Efna5 | Ephb2 | ECs
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Efna5"
$ws.Range("C2").Value2 = "Ephb2"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.05800433333333333
$ws.Range("H2").Value2 = 0.174013
$ws.Range("I2").Value2 = 0.02087975181349295
$ws.Range("J2").Value2 = 0.02087975181349295
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.0006136666666666667
$ws.Range("N2").Value2 = 0.001841
$ws.Range("O2").Value2 = 0.000129696697123199
$ws.Range("P2").Value2 = 0.000129696697123199
$ws.Range("Q2").Value2 = 0.00003559532588888889
$ws.Range("R2").Value2 = 0.000320357933
$ws.Range("S2").Value2 = 0.000002708034846962161
$ws.Range("T2").Value2 = 0.000002708034846962161

# Row 3: ECs | Efna5 | Ephb2 | FAPs
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Efna5"
$ws.Range("C3").Value2 = "Ephb2"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.05800433333333333
$ws.Range("H3").Value2 = 0.174013
$ws.Range("I3").Value2 = 0.02087975181349295
$ws.Range("J3").Value2 = 0.02087975181349295
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 3.821776
$ws.Range("N3").Value2 = 11.465328
$ws.Range("O3").Value2 = 0.8077214410831794
$ws.Range("P3").Value2 = 0.8077214410831794
$ws.Range("Q3").Value2 = 0.2216795690293333
$ws.Range("R3").Value2 = 1.995116121264
$ws.Range("S3").Value2 = 0.01686502322425366
$ws.Range("T3").Value2 = 0.01686502322425366

# Row 4: ECs | Efna5 | Ephb2 | MuSCs
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Efna5"
$ws.Range("C4").Value2 = "Ephb2"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.05800433333333333
$ws.Range("H4").Value2 = 0.174013
$ws.Range("I4").Value2 = 0.02087975181349295
$ws.Range("J4").Value2 = 0.02087975181349295
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 0.9091623333333333
$ws.Range("N4").Value2 = 2.727487
$ws.Range("O4").Value2 = 0.1921488622196973
$ws.Range("P4").Value2 = 0.1921488622196973
$ws.Range("Q4").Value2 = 0.05273535503677777
$ws.Range("R4").Value2 = 0.474618195331
$ws.Range("S4").Value2 = 0.004012020554392332
$ws.Range("T4").Value2 = 0.004012020554392332

# Row 5: FAPs | Efna5 | Ephb2 | ECs
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Efna5"
$ws.Range("C5").Value2 = "Ephb2"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 1.666083666666667
$ws.Range("H5").Value2 = 4.998251
$ws.Range("I5").Value2 = 0.5997381826733804
$ws.Range("J5").Value2 = 0.5997381826733805
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.0006136666666666667
$ws.Range("N5").Value2 = 0.001841
$ws.Range("O5").Value2 = 0.000129696697123199
$ws.Range("P5").Value2 = 0.000129696697123199
$ws.Range("Q5").Value2 = 0.001022420010111111
$ws.Range("R5").Value2 = 0.009201780090999999
$ws.Range("S5").Value2 = 0.00007778406143140724
$ws.Range("T5").Value2 = 0.00007778406143140726

# Row 6: FAPs | Efna5 | Ephb2 | FAPs
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Efna5"
$ws.Range("C6").Value2 = "Ephb2"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 1.666083666666667
$ws.Range("H6").Value2 = 4.998251
$ws.Range("I6").Value2 = 0.5997381826733804
$ws.Range("J6").Value2 = 0.5997381826733805
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 3.821776
$ws.Range("N6").Value2 = 11.465328
$ws.Range("O6").Value2 = 0.8077214410831794
$ws.Range("P6").Value2 = 0.8077214410831794
$ws.Range("Q6").Value2 = 6.367398571258666
$ws.Range("R6").Value2 = 57.306587141328
$ws.Range("S6").Value2 = 0.48442138918155
$ws.Range("T6").Value2 = 0.48442138918155

# Row 7: FAPs | Efna5 | Ephb2 | MuSCs
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Efna5"
$ws.Range("C7").Value2 = "Ephb2"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 1.666083666666667
$ws.Range("H7").Value2 = 4.998251
$ws.Range("I7").Value2 = 0.5997381826733804
$ws.Range("J7").Value2 = 0.5997381826733805
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 0.9091623333333333
$ws.Range("N7").Value2 = 2.727487
$ws.Range("O7").Value2 = 0.1921488622196973
$ws.Range("P7").Value2 = 0.1921488622196973
$ws.Range("Q7").Value2 = 1.514740513915222
$ws.Range("R7").Value2 = 13.632664625237
$ws.Range("S7").Value2 = 0.115239009430399
$ws.Range("T7").Value2 = 0.1152390094303991

# Row 8: MuSCs | Efna5 | Ephb2 | ECs
$ws.Range("A8").Value2 = "MuSCs"
$ws.Range("B8").Value2 = "Efna5"
$ws.Range("C8").Value2 = "Ephb2"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 1.053930333333333
$ws.Range("H8").Value2 = 3.161791
$ws.Range("I8").Value2 = 0.3793820655131266
$ws.Range("J8").Value2 = 0.3793820655131266
$ws.Range("K8").Value2 = 1
$ws.Range("L8").Value2 = 0.3333333333333333
$ws.Range("M8").Value2 = 0.0006136666666666667
$ws.Range("N8").Value2 = 0.001841
$ws.Range("O8").Value2 = 0.000129696697123199
$ws.Range("P8").Value2 = 0.000129696697123199
$ws.Range("Q8").Value2 = 0.0006467619145555556
$ws.Range("R8").Value2 = 0.005820857231
$ws.Range("S8").Value2 = 0.00004920460084482963
$ws.Range("T8").Value2 = 0.00004920460084482963

# Row 9: MuSCs | Efna5 | Ephb2 | FAPs
$ws.Range("A9").Value2 = "MuSCs"
$ws.Range("B9").Value2 = "Efna5"
$ws.Range("C9").Value2 = "Ephb2"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 1.053930333333333
$ws.Range("H9").Value2 = 3.161791
$ws.Range("I9").Value2 = 0.3793820655131266
$ws.Range("J9").Value2 = 0.3793820655131266
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 3.821776
$ws.Range("N9").Value2 = 11.465328
$ws.Range("O9").Value2 = 0.8077214410831794
$ws.Range("P9").Value2 = 0.8077214410831794
$ws.Range("Q9").Value2 = 4.027885653605334
$ws.Range("R9").Value2 = 36.250970882448
$ws.Range("S9").Value2 = 0.3064350286773758
$ws.Range("T9").Value2 = 0.3064350286773758

# Row 10: MuSCs | Efna5 | Ephb2 | MuSCs
$ws.Range("A10").Value2 = "MuSCs"
$ws.Range("B10").Value2 = "Efna5"
$ws.Range("C10").Value2 = "Ephb2"
$ws.Range("D10").Value2 = "MuSCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 1.053930333333333
$ws.Range("H10").Value2 = 3.161791
$ws.Range("I10").Value2 = 0.3793820655131266
$ws.Range("J10").Value2 = 0.3793820655131266
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.9091623333333333
$ws.Range("N10").Value2 = 2.727487
$ws.Range("O10").Value2 = 0.1921488622196973
$ws.Range("P10").Value2 = 0.1921488622196973
$ws.Range("Q10").Value2 = 0.9581937610241111
$ws.Range("R10").Value2 = 8.623743849217
$ws.Range("S10").Value2 = 0.07289783223490594
$ws.Range("T10").Value2 = 0.07289783223490594

"Updated Efna5-Ephb2 NATMI data: rows 2-10 now populated (added ECs sender rows, refreshed FAPs/MuSCs rows)."